$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 72 -----------------------------------------------------------
# This row currently carries leftover empty placeholder cells (I72 and
# K72:R72) from an earlier edit. The new revision drops them so the row
# only spans through column J, matching every other row in the sheet.
$ws.Range("I72").ClearContents()
$ws.Range("K72:R72").ClearContents()

# --- Row 73 -------------------------------------------------------------
# New "Solo Revisión" entry, columns A-H + J only (no extra cells).
$ws.Range("A73").Value = "2CA06708"
$ws.Range("B73").Value = "REVLON COLORSTAY MASC.COLOR BOOSTER BLONDE 125ML"
$ws.Range("C73").Value = "CABELLO ACONDIC. SUAVIZANTE"
$ws.Range("D73").Value = "Tiene PT"
$ws.Range("E73").Value = "Tiene ES"
$ws.Range("F73").Value = "No Tiene IT - TRADOTTO"
$ws.Range("G73").Value = "'125"
$ws.Range("H73").Value = "ML"
$ws.Range("J73").Value = "Solo Revisión"

# --- Row 74 -------------------------------------------------------------
# New "Revisado y Traducido" entry, columns A-H + J populated, and
# I74 / K74:R74 kept as present-but-empty placeholder cells (mirrors the
# pattern row 72 used to have before this revision).
$ws.Range("A74").Value = "2CN01618"
$ws.Range("B74").Value = "CREME OF NATURE HONEY CREMA DEFINICION RIZOS 326GR"
$ws.Range("C74").Value = "CABELLO TONICO LOCION"
$ws.Range("D74").Value = "No Tiene PT - TRADUZIDO"
$ws.Range("E74").Value = "Tiene ES"
$ws.Range("F74").Value = "No Tiene IT - TRADOTTO"
$ws.Range("G74").Value = "'326"
$ws.Range("H74").Value = "GR"
$ws.Range("I74").Value = "'"
$ws.Range("J74").Value = "Revisado y Traducido"
$ws.Range("K74:R74").Value = "'"

# The quantity cells (G73/G74) and the blank placeholder cells (I74,
# K74:R74) were written with a leading apostrophe so they are stored as
# text (matching the rest of the sheet, which is text-typed throughout)
# instead of being auto-converted to numbers / true blanks. That leaves
# a stray "quote prefix" look on the cells themselves, so restore their
# display format from an untouched, plain cell without touching the
# values that were just entered.
$ws.Range("A1").Copy()
$ws.Range("G73").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Copy()
$ws.Range("G74").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Copy()
$ws.Range("I74").PasteSpecial(-4122)   # xlPasteFormats
$ws.Range("A1").Copy()
$ws.Range("K74:R74").PasteSpecial(-4122)   # xlPasteFormats
